$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 78977
$ws.Range("E2").Value = 5934
$ws.Range("F2").Value = 5934
$ws.Range("G2").Value = 6368
$ws.Range("H2").Value = 4343
$ws.Range("I2").Value = 4129
$ws.Range("J2").Value = 214
$ws.Range("K2").Value = 55460
$ws.Range("L2").Value = 13320
$ws.Range("M2").Value = 42140
$ws.Range("N2").Value = 41123
$ws.Range("O2").Value = 1018
$ws.Range("P2").Value = 387
$ws.Range("Q2").Value = 9963
$ws.Range("R2").Value = -8403
$ws.Range("S2").Value = -510
$ws.Range("T2").Value = 3839
$ws.Range("U2").Value = 6124
$ws.Range("V2").Value = 202
$ws.Range("W2").Value = 7.51
$ws.Range("X2").Value = 5.5
$ws.Range("Y2").Value = 10.43
$ws.Range("Z2").Value = 8.01
$ws.Range("AA2").Value = 31.61
$ws.Range("AB2").Value = 11237.25
$ws.Range("AC2").Value = 5336
$ws.Range("AD2").Value = 55
$ws.Range("AE2").Value = 53164
$ws.Range("AF2").Value = 5.52
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 0.17
$ws.Range("AI2").Value = 9.369999999999999
$ws.Range("AJ2").Value = 77377800
$ws.Range("D3").Value = 78535
$ws.Range("E3").Value = 5883
$ws.Range("F3").Value = 5883
$ws.Range("G3").Value = 6491
$ws.Range("H3").Value = 4698
$ws.Range("I3").Value = 4390
$ws.Range("J3").Value = 308
$ws.Range("K3").Value = 63317
$ws.Range("L3").Value = 15548
$ws.Range("M3").Value = 47769
$ws.Range("N3").Value = 45888
$ws.Range("O3").Value = 1881
$ws.Range("P3").Value = 387
$ws.Range("Q3").Value = 8648
$ws.Range("R3").Value = -7824
$ws.Range("S3").Value = -481
$ws.Range("T3").Value = 3618
$ws.Range("U3").Value = 5029
$ws.Range("V3").Value = 108
$ws.Range("W3").Value = 7.49
$ws.Range("X3").Value = 5.98
$ws.Range("Y3").Value = 10.09
$ws.Range("Z3").Value = 7.91
$ws.Range("AA3").Value = 32.55
$ws.Range("AB3").Value = 12269.66
$ws.Range("AC3").Value = 5674
$ws.Range("AD3").Value = 44.77
$ws.Range("AE3").Value = 59324
$ws.Range("AF3").Value = 4.28
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 0.2
$ws.Range("AI3").Value = 8.81
$ws.Range("AJ3").Value = 77377800
$ws.Range("D4").Value = 81802
$ws.Range("E4").Value = 6271
$ws.Range("F4").Value = 6271
$ws.Range("G4").Value = 7519
$ws.Range("H4").Value = 5143
$ws.Range("I4").Value = 4639
$ws.Range("J4").Value = 505
$ws.Range("K4").Value = 68420
$ws.Range("L4").Value = 15509
$ws.Range("M4").Value = 52911
$ws.Range("N4").Value = 51407
$ws.Range("O4").Value = 1504
$ws.Range("P4").Value = 387
$ws.Range("Q4").Value = 6228
$ws.Range("R4").Value = -2539
$ws.Range("S4").Value = -443
$ws.Range("T4").Value = 1898
$ws.Range("U4").Value = 4330
$ws.Range("V4").Value = 67
$ws.Range("W4").Value = 7.67
$ws.Range("X4").Value = 6.29
$ws.Range("Y4").Value = 9.539999999999999
$ws.Range("Z4").Value = 7.81
$ws.Range("AA4").Value = 29.31
$ws.Range("AB4").Value = 13368.64
$ws.Range("AC4").Value = 5995
$ws.Range("AD4").Value = 23.27
$ws.Range("AE4").Value = 66460
$ws.Range("AF4").Value = 2.1
$ws.Range("AG4").Value = 750
$ws.Range("AH4").Value = 0.54
$ws.Range("AI4").Value = 12.51
$ws.Range("AJ4").Value = 77377800
$ws.Range("D5").Value = 92992
$ws.Range("E5").Value = 7316
$ws.Range("F5").Value = 7316
$ws.Range("G5").Value = 7521
$ws.Range("H5").Value = 5418
$ws.Range("I5").Value = 5303
$ws.Range("J5").Value = 114
$ws.Range("K5").Value = 72778
$ws.Range("L5").Value = 15583
$ws.Range("M5").Value = 57194
$ws.Range("N5").Value = 55583
$ws.Range("O5").Value = 1612
$ws.Range("P5").Value = 387
$ws.Range("Q5").Value = 8348
$ws.Range("R5").Value = -9704
$ws.Range("S5").Value = -676
$ws.Range("T5").Value = 1921
$ws.Range("U5").Value = 6428
$ws.Range("V5").Value = 8
$ws.Range("W5").Value = 7.87
$ws.Range("X5").Value = 5.83
$ws.Range("Y5").Value = 9.91
$ws.Range("Z5").Value = 7.67
$ws.Range("AA5").Value = 27.25
$ws.Range("AB5").Value = 14589.49
$ws.Range("AC5").Value = 6854
$ws.Range("AD5").Value = 29.18
$ws.Range("AE5").Value = 71859
$ws.Range("AF5").Value = 2.78
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1
$ws.Range("AI5").Value = 29.17
$ws.Range("AJ5").Value = 77377800
$ws.Range("D6").Value = 100342
$ws.Range("E6").Value = 8774
$ws.Range("F6").Value = 8774
$ws.Range("G6").Value = 9565
$ws.Range("H6").Value = 6388
$ws.Range("I6").Value = 6294
$ws.Range("K6").Value = 80138
$ws.Range("L6").Value = 18653
$ws.Range("M6").Value = 61486
$ws.Range("N6").Value = 59825
$ws.Range("P6").Value = 387
$ws.Range("Q6").Value = 11965
$ws.Range("R6").Value = -8186
$ws.Range("S6").Value = -1562
$ws.Range("T6").Value = 2514
$ws.Range("U6").Value = 9451
$ws.Range("V6").Value = 8
$ws.Range("W6").Value = 8.74
$ws.Range("X6").Value = 6.37
$ws.Range("Y6").Value = 10.91
$ws.Range("Z6").Value = 8.359999999999999
$ws.Range("AA6").Value = 30.34
$ws.Range("AB6").Value = 15816.46
$ws.Range("AC6").Value = 8134
$ws.Range("AD6").Value = 25.08
$ws.Range("AE6").Value = 77342
$ws.Range("AF6").Value = 2.64
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 0.98
$ws.Range("AI6").Value = 24.58
$ws.Range("AJ6").Value = 77377800
$ws.Range("D7").Value = 109463
$ws.Range("E7").Value = 9638
$ws.Range("G7").Value = 10408
$ws.Range("H7").Value = 7140
$ws.Range("I7").Value = 7022
$ws.Range("K7").Value = 90474
$ws.Range("L7").Value = 23272
$ws.Range("M7").Value = 67203
$ws.Range("N7").Value = 65432
$ws.Range("P7").Value = 389
$ws.Range("Q7").Value = 11174
$ws.Range("R7").Value = -4368
$ws.Range("S7").Value = -2034
$ws.Range("T7").Value = 2880
$ws.Range("U7").Value = 7692
$ws.Range("W7").Value = 8.800000000000001
$ws.Range("X7").Value = 6.52
$ws.Range("Y7").Value = 11.21
$ws.Range("Z7").Value = 8.369999999999999
$ws.Range("AA7").Value = 34.63
$ws.Range("AC7").Value = 9075
$ws.Range("AD7").Value = 22.59
$ws.Range("AE7").Value = 84592
$ws.Range("AF7").Value = 2.42
$ws.Range("AG7").Value = 2110
$ws.Range("AH7").Value = 1.03
$ws.Range("AI7").Value = 23.25
$ws.Range("D8").Value = 117267
$ws.Range("E8").Value = 10708
$ws.Range("G8").Value = 11457
$ws.Range("H8").Value = 8065
$ws.Range("I8").Value = 7947
$ws.Range("K8").Value = 98775
$ws.Range("L8").Value = 24953
$ws.Range("M8").Value = 73822
$ws.Range("N8").Value = 71933
$ws.Range("P8").Value = 389
$ws.Range("Q8").Value = 11882
$ws.Range("R8").Value = -4702
$ws.Range("S8").Value = -1743
$ws.Range("T8").Value = 2159
$ws.Range("U8").Value = 8748
$ws.Range("W8").Value = 9.130000000000001
$ws.Range("X8").Value = 6.88
$ws.Range("Y8").Value = 11.59
$ws.Range("Z8").Value = 8.56
$ws.Range("AA8").Value = 33.8
$ws.Range("AC8").Value = 10270
$ws.Range("AD8").Value = 18.89
$ws.Range("AE8").Value = 92996
$ws.Range("AF8").Value = 2.09
$ws.Range("AG8").Value = 2380
$ws.Range("AH8").Value = 1.23
$ws.Range("AI8").Value = 23.17
$ws.Range("D9").Value = 127626
$ws.Range("E9").Value = 11739
$ws.Range("G9").Value = 12527
$ws.Range("H9").Value = 8845
$ws.Range("I9").Value = 8693
$ws.Range("K9").Value = 107511
$ws.Range("L9").Value = 26727
$ws.Range("M9").Value = 80784
$ws.Range("N9").Value = 78704
$ws.Range("P9").Value = 389
$ws.Range("Q9").Value = 12195
$ws.Range("R9").Value = -5116
$ws.Range("S9").Value = -1844
$ws.Range("T9").Value = 2105
$ws.Range("U9").Value = 10331
$ws.Range("W9").Value = 9.199999999999999
$ws.Range("X9").Value = 6.93
$ws.Range("Y9").Value = 11.54
$ws.Range("Z9").Value = 8.58
$ws.Range("AA9").Value = 33.08
$ws.Range("AC9").Value = 11235
$ws.Range("AD9").Value = 17.27
$ws.Range("AE9").Value = 101750
$ws.Range("AF9").Value = 1.91
$ws.Range("AG9").Value = 2558
$ws.Range("AH9").Value = 1.32
$ws.Range("AI9").Value = 22.77
